$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9812178611755371
$ws.Range("B1").Value = 1.122846484184265
$ws.Range("C1").Value = 5.326637744903564
$ws.Range("D1").Value = 1.58436119556427
$ws.Range("E1").Value = 0.9647256135940552
